$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Widen column A slightly (ColumnWidth -> stored xlsx width is offset by the
# default font metrics; 20.14 lands on a stored width of exactly 21)
$ws.Columns.Item(1).ColumnWidth = 20.14

# Header row: name and week number
$ws.Range("C1").Value = "Richard Dobson"
$ws.Range("E1").Value = 5

# Task rows 3-6
$ws.Range("A3").Value = "Project Design and Specifics"
$ws.Range("B3").Value = "Design Command Line UI"
$ws.Range("C3").Value = 2
$ws.Range("D3").Value = 2

$ws.Range("A4").Value = "Project Design and Specifics"
$ws.Range("B4").Value = "Analysis of Open source Metadata Libraries and Programs"
$ws.Range("C4").Value = 4
$ws.Range("D4").Value = 4

$ws.Range("A5").Value = "Project Design and Specifics"
$ws.Range("B5").Value = "Familiarisation with libraries (native and external)"
$ws.Range("C5").Value = 6
$ws.Range("D5").Value = 6

$ws.Range("A6").Value = "Project Build"
$ws.Range("B6").Value = "Learn Python Syntax"
$ws.Range("C6").Value = 8
$ws.Range("D6").Value = 8

# Apply the same "wrap text" style used by row 3 to row 7's Stage cell, matching the others,
# even though row 7 has no content yet.
$ws.Range("A3:A7").WrapText = $true

# Rows 3-5 now hold two-line task descriptions, so they grow taller to fit
$ws.Rows.Item(3).RowHeight = 27.75
$ws.Rows.Item(4).RowHeight = 27.75
$ws.Rows.Item(5).RowHeight = 27.75

# Update cumulative total label
$ws.Range("A14").Value = "Cumulative Total: 100"

# Update the active selection to match the author's last position
$ws.Range("B28").Select()
